$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Get-ParaByExactText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# Escape text for safe embedding inside XML (handles & < > etc.)
function Xml-Escape($s) {
    $s = $s -replace "&", "&amp;"
    $s = $s -replace "<", "&lt;"
    $s = $s -replace ">", "&gt;"
    return $s
}

$newTitle = "Play Jungle Giants Free - Unique Gameplay & Stunning Graphics"
$newTitleEsc = Xml-Escape $newTitle

# 1. Main heading (Heading1), no preceding empty run -> plain Find/Replace is safe.
$d.Content.Find.Execute("Play Jungle Giants Free - A Unique Online Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, $newTitle, 2)

# 2. Pros bullet: "Immersive sound design" -> "Impressive sound design"
#    Paragraph has a leading empty <w:r/> that must be preserved, so rebuild the
#    whole paragraph via InsertXML instead of a plain text replace (which would
#    silently drop the empty run).
$p = Get-ParaByExactText $d "Immersive sound design"
if ($p -ne $null) {
    $xml = "<w:p $wns><w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Impressive sound design</w:t></w:r></w:p>"
    $p.Range.InsertXML($xml)
}

# 3. Cons bullet: "Limited range of animal symbols" -> "Limited number of free spins"
$p = Get-ParaByExactText $d "Limited range of animal symbols"
if ($p -ne $null) {
    $xml = "<w:p $wns><w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Limited number of free spins</w:t></w:r></w:p>"
    $p.Range.InsertXML($xml)
}

# 4. Bold meta-title paragraph (same text as the H1 heading)
$p = Get-ParaByExactText $d "Play Jungle Giants Free - A Unique Online Slot Game"
if ($p -ne $null) {
    $xml = "<w:p $wns><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>$newTitleEsc</w:t></w:r></w:p>"
    $p.Range.InsertXML($xml)
}

# 5. Italic meta-description paragraph
$p = Get-ParaByExactText $d "Read our review of Jungle Giants, a unique online slot game with stunning graphics, immersive sound, and numerous winning opportunities. Play free now."
if ($p -ne $null) {
    $newDesc = "Read our review of Jungle Giants and play this unique online slot for free. Experience stunning graphics and numerous winning opportunities."
    $newDescEsc = Xml-Escape $newDesc
    $xml = "<w:p $wns><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>$newDescEsc</w:t></w:r></w:p>"
    $p.Range.InsertXML($xml)
}

"done"
